$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of quotations for 2025-09-21 (Excel serial date 45921)
$ws.Range("A17").Value = 45921
$ws.Range("A17").NumberFormat = $ws.Range("A16").NumberFormat

$ws.Range("B17").Value = "20,9437"
$ws.Range("C17").Value = "15,0727"
$ws.Range("D17").Value = "14,9476"
$ws.Range("E17").Value = "14,9476"
